$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 7).Value = 38.552387
$ws.Cells.Item(2, 8).Value = 115.657161
$ws.Cells.Item(2, 9).Value = 0.6603131536923842
$ws.Cells.Item(2, 10).Value = 0.6603131536923841
$ws.Cells.Item(2, 13).Value = 41.46003700000001
$ws.Cells.Item(2, 14).Value = 124.380111
$ws.Cells.Item(2, 15).Value = 0.5916248119519026
$ws.Cells.Item(2, 16).Value = 0.5916248119519028
$ws.Cells.Item(2, 17).Value = 1598.383391458319
$ws.Cells.Item(2, 18).Value = 14385.45052312487
$ws.Cells.Item(2, 19).Value = 0.3906576453826246
$ws.Cells.Item(2, 20).Value = 0.3906576453826246
$ws.Cells.Item(3, 7).Value = 38.552387
$ws.Cells.Item(3, 8).Value = 115.657161
$ws.Cells.Item(3, 9).Value = 0.6603131536923842
$ws.Cells.Item(3, 10).Value = 0.6603131536923841
$ws.Cells.Item(3, 15).Value = 0.01390494488024241
$ws.Cells.Item(3, 16).Value = 0.01390494488024241
$ws.Cells.Item(3, 17).Value = 37.56676952475366
$ws.Cells.Item(3, 18).Value = 338.100925722783
$ws.Cells.Item(3, 19).Value = 0.009181618005791635
$ws.Cells.Item(3, 20).Value = 0.009181618005791637
$ws.Cells.Item(4, 7).Value = 38.552387
$ws.Cells.Item(4, 8).Value = 115.657161
$ws.Cells.Item(4, 9).Value = 0.6603131536923842
$ws.Cells.Item(4, 10).Value = 0.6603131536923841
$ws.Cells.Item(4, 13).Value = 27.57046566666667
$ws.Cells.Item(4, 14).Value = 82.71139700000001
$ws.Cells.Item(4, 15).Value = 0.3934239510077632
$ws.Cells.Item(4, 16).Value = 0.3934239510077633
$ws.Cells.Item(4, 17).Value = 1062.907262151546
$ws.Cells.Item(4, 18).Value = 9566.165359363917
$ws.Cells.Item(4, 19).Value = 0.2597830098280542
$ws.Cells.Item(4, 20).Value = 0.2597830098280542
$ws.Cells.Item(5, 7).Value = 38.552387
$ws.Cells.Item(5, 8).Value = 115.657161
$ws.Cells.Item(5, 9).Value = 0.6603131536923842
$ws.Cells.Item(5, 10).Value = 0.6603131536923841
$ws.Cells.Item(5, 11).Value = 3
$ws.Cells.Item(5, 12).Value = 1
$ws.Cells.Item(5, 13).Value = 0.07332233333333334
$ws.Cells.Item(5, 14).Value = 0.219967
$ws.Cells.Item(5, 15).Value = 0.001046292160091609
$ws.Cells.Item(5, 16).Value = 0.00104629216009161
$ws.Cells.Item(5, 17).Value = 2.826750970409666
$ws.Cells.Item(5, 18).Value = 25.440758733687
$ws.Cells.Item(5, 19).Value = 0.0006908804759137075
$ws.Cells.Item(5, 20).Value = 0.0006908804759137075
$ws.Cells.Item(6, 9).Value = 0.012622453244502
$ws.Cells.Item(6, 10).Value = 0.012622453244502
$ws.Cells.Item(6, 13).Value = 41.46003700000001
$ws.Cells.Item(6, 14).Value = 124.380111
$ws.Cells.Item(6, 15).Value = 0.5916248119519026
$ws.Cells.Item(6, 16).Value = 0.5916248119519028
$ws.Cells.Item(6, 17).Value = 30.554471787594
$ws.Cells.Item(6, 18).Value = 274.990246088346
$ws.Cells.Item(6, 19).Value = 0.007467756527150182
$ws.Cells.Item(6, 20).Value = 0.007467756527150182
$ws.Cells.Item(7, 9).Value = 0.012622453244502
$ws.Cells.Item(7, 10).Value = 0.012622453244502
$ws.Cells.Item(7, 15).Value = 0.01390494488024241
$ws.Cells.Item(7, 16).Value = 0.01390494488024241
$ws.Cells.Item(7, 19).Value = 0.0001755145166182373
$ws.Cells.Item(7, 20).Value = 0.0001755145166182373
$ws.Cells.Item(8, 9).Value = 0.012622453244502
$ws.Cells.Item(8, 10).Value = 0.012622453244502
$ws.Cells.Item(8, 13).Value = 27.57046566666667
$ws.Cells.Item(8, 14).Value = 82.71139700000001
$ws.Cells.Item(8, 15).Value = 0.3934239510077632
$ws.Cells.Item(8, 16).Value = 0.3934239510077633
$ws.Cells.Item(8, 17).Value = 20.318385518638
$ws.Cells.Item(8, 18).Value = 182.865469667742
$ws.Cells.Item(8, 19).Value = 0.004965975426862739
$ws.Cells.Item(8, 20).Value = 0.004965975426862739
$ws.Cells.Item(9, 9).Value = 0.012622453244502
$ws.Cells.Item(9, 10).Value = 0.012622453244502
$ws.Cells.Item(9, 11).Value = 3
$ws.Cells.Item(9, 12).Value = 1
$ws.Cells.Item(9, 13).Value = 0.07332233333333334
$ws.Cells.Item(9, 14).Value = 0.219967
$ws.Cells.Item(9, 15).Value = 0.001046292160091609
$ws.Cells.Item(9, 16).Value = 0.00104629216009161
$ws.Cells.Item(9, 17).Value = 0.054035773418
$ws.Cells.Item(9, 18).Value = 0.486321960762
$ws.Cells.Item(9, 19).Value = 0.00001320677387084535
$ws.Cells.Item(9, 20).Value = 0.00001320677387084535
$ws.Cells.Item(10, 7).Value = 18.94833433333333
$ws.Cells.Item(10, 8).Value = 56.845003
$ws.Cells.Item(10, 9).Value = 0.3245411082032615
$ws.Cells.Item(10, 10).Value = 0.3245411082032615
$ws.Cells.Item(10, 13).Value = 41.46003700000001
$ws.Cells.Item(10, 14).Value = 124.380111
$ws.Cells.Item(10, 15).Value = 0.5916248119519026
$ws.Cells.Item(10, 16).Value = 0.5916248119519028
$ws.Cells.Item(10, 17).Value = 785.5986425483704
$ws.Cells.Item(10, 18).Value = 7070.387782935333
$ws.Cells.Item(10, 19).Value = 0.1920065721114167
$ws.Cells.Item(10, 20).Value = 0.1920065721114167
$ws.Cells.Item(11, 7).Value = 18.94833433333333
$ws.Cells.Item(11, 8).Value = 56.845003
$ws.Cells.Item(11, 9).Value = 0.3245411082032615
$ws.Cells.Item(11, 10).Value = 0.3245411082032615
$ws.Cells.Item(11, 15).Value = 0.01390494488024241
$ws.Cells.Item(11, 16).Value = 0.01390494488024241
$ws.Cells.Item(11, 17).Value = 18.46390753387877
$ws.Cells.Item(11, 18).Value = 166.175167804909
$ws.Cells.Item(11, 19).Value = 0.004512726220939138
$ws.Cells.Item(11, 20).Value = 0.00451272622093914
$ws.Cells.Item(12, 7).Value = 18.94833433333333
$ws.Cells.Item(12, 8).Value = 56.845003
$ws.Cells.Item(12, 9).Value = 0.3245411082032615
$ws.Cells.Item(12, 10).Value = 0.3245411082032615
$ws.Cells.Item(12, 13).Value = 27.57046566666667
$ws.Cells.Item(12, 14).Value = 82.71139700000001
$ws.Cells.Item(12, 15).Value = 0.3934239510077632
$ws.Cells.Item(12, 16).Value = 0.3934239510077633
$ws.Cells.Item(12, 17).Value = 522.4144011776879
$ws.Cells.Item(12, 18).Value = 4701.729610599191
$ws.Cells.Item(12, 19).Value = 0.1276822450537652
$ws.Cells.Item(12, 20).Value = 0.1276822450537652
$ws.Cells.Item(13, 7).Value = 18.94833433333333
$ws.Cells.Item(13, 8).Value = 56.845003
$ws.Cells.Item(13, 9).Value = 0.3245411082032615
$ws.Cells.Item(13, 10).Value = 0.3245411082032615
$ws.Cells.Item(13, 11).Value = 3
$ws.Cells.Item(13, 12).Value = 1
$ws.Cells.Item(13, 13).Value = 0.07332233333333334
$ws.Cells.Item(13, 14).Value = 0.219967
$ws.Cells.Item(13, 15).Value = 0.001046292160091609
$ws.Cells.Item(13, 16).Value = 0.00104629216009161
$ws.Cells.Item(13, 17).Value = 1.389336086100111
$ws.Cells.Item(13, 18).Value = 12.504024774901
$ws.Cells.Item(13, 19).Value = 0.0003395648171405152
$ws.Cells.Item(13, 20).Value = 0.0003395648171405153
$ws.Cells.Item(14, 7).Value = 0.147322
$ws.Cells.Item(14, 8).Value = 0.441966
$ws.Cells.Item(14, 9).Value = 0.002523284859852372
$ws.Cells.Item(14, 10).Value = 0.002523284859852372
$ws.Cells.Item(14, 13).Value = 41.46003700000001
$ws.Cells.Item(14, 14).Value = 124.380111
$ws.Cells.Item(14, 15).Value = 0.5916248119519026
$ws.Cells.Item(14, 16).Value = 0.5916248119519028
$ws.Cells.Item(14, 17).Value = 6.107975570914
$ws.Cells.Item(14, 18).Value = 54.971780138226
$ws.Cells.Item(14, 19).Value = 0.001492837930711243
$ws.Cells.Item(14, 20).Value = 0.001492837930711243
$ws.Cells.Item(15, 7).Value = 0.147322
$ws.Cells.Item(15, 8).Value = 0.441966
$ws.Cells.Item(15, 9).Value = 0.002523284859852372
$ws.Cells.Item(15, 10).Value = 0.002523284859852372
$ws.Cells.Item(15, 15).Value = 0.01390494488024241
$ws.Cells.Item(15, 16).Value = 0.01390494488024241
$ws.Cells.Item(15, 17).Value = 0.1435556148553333
$ws.Cells.Item(15, 18).Value = 1.292000533698
$ws.Cells.Item(15, 19).Value = 0.00003508613689339742
$ws.Cells.Item(15, 20).Value = 0.00003508613689339743
$ws.Cells.Item(16, 7).Value = 0.147322
$ws.Cells.Item(16, 8).Value = 0.441966
$ws.Cells.Item(16, 9).Value = 0.002523284859852372
$ws.Cells.Item(16, 10).Value = 0.002523284859852372
$ws.Cells.Item(16, 13).Value = 27.57046566666667
$ws.Cells.Item(16, 14).Value = 82.71139700000001
$ws.Cells.Item(16, 15).Value = 0.3934239510077632
$ws.Cells.Item(16, 16).Value = 0.3934239510077633
$ws.Cells.Item(16, 17).Value = 4.061736142944667
$ws.Cells.Item(16, 18).Value = 36.555625286502
$ws.Cells.Item(16, 19).Value = 0.0009927206990811904
$ws.Cells.Item(16, 20).Value = 0.0009927206990811906
$ws.Cells.Item(17, 7).Value = 0.147322
$ws.Cells.Item(17, 8).Value = 0.441966
$ws.Cells.Item(17, 9).Value = 0.002523284859852372
$ws.Cells.Item(17, 10).Value = 0.002523284859852372
$ws.Cells.Item(17, 11).Value = 3
$ws.Cells.Item(17, 12).Value = 1
$ws.Cells.Item(17, 13).Value = 0.07332233333333334
$ws.Cells.Item(17, 14).Value = 0.219967
$ws.Cells.Item(17, 15).Value = 0.001046292160091609
$ws.Cells.Item(17, 16).Value = 0.00104629216009161
$ws.Cells.Item(17, 17).Value = 0.01080199279133333
$ws.Cells.Item(17, 18).Value = 0.09721793512199998
$ws.Cells.Item(17, 19).Value = 0.000002640093166541392
$ws.Cells.Item(17, 20).Value = 0.000002640093166541393
